# Updates loading_percent values for rows 2-25 (columns B,C,D,F,G,I,J,K,L,M,O)
# to reflect the re-run "case with 380 kV done" results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","I","J","K","L","M","O")

$data = @(
    @(2, 12.4792979961822, 7.93558209730972, 8.424016779254796, 38.43214631172326, 3.696992972455111, 29.06301824678741, 10.87162690171373, 9.746763163101811, 11.72846626682434, 15.29062946722657, 29.47963852309263),
    @(3, 12.27646224116472, 7.896550223921045, 8.413461695005585, 38.53557918211512, 3.698740068981173, 29.16373847256676, 10.89162479851963, 9.597112243137619, 11.73685527144248, 15.26288465786565, 29.57365006101284),
    @(4, 12.15207091575254, 7.872211261532594, 8.408019066720966, 38.60608122061988, 3.699870451771593, 29.22996718555612, 10.90458675209637, 9.505327402630574, 11.74328114455768, 15.24762871338478, 29.63626498576981),
    @(5, 12.10148066502587, 7.862201868319188, 8.406064136626254, 38.63656859509847, 3.700345634794531, 29.25805935554513, 10.9100411076358, 9.467993466641504, 11.74622108367374, 15.24186364669993, 29.66301067275804),
    @(6, 12.09308807880952, 7.860534397080345, 8.40575545899331, 38.64173708762859, 3.700425418226363, 29.26279069902854, 10.91095721670311, 9.461799651405336, 11.74672868738742, 15.24093377586143, 29.66752602638395),
    @(7, 12.15138814933827, 7.872076636414854, 8.407991634594387, 38.60648527109554, 3.699876801313493, 29.23034157794061, 10.90465961336959, 9.504823566330932, 11.74331949155088, 15.24754912852918, 29.63662070941387),
    @(8, 12.40936615013307, 7.922203039784821, 8.420163083321368, 38.46635757135143, 3.697583430953381, 29.09683655764572, 10.8783806673635, 9.695168527063158, 11.73109462480651, 15.28069643172624, 29.51103824517437),
    @(9, 12.91375953415896, 8.017418903157793, 8.452180757794739, 38.24712070963571, 3.693541616999373, 28.86981189216834, 10.83224703277197, 10.0673478929454, 11.71720494334051, 15.35961562358296, 29.30360031161193),
    @(10, 13.27990283222834, 8.085336744759342, 8.480547284031715, 38.11999009142019, 3.690846946503629, 28.72418612050431, 10.80161440744329, 10.33764874572119, 11.71310126770668, 15.42579819721248, 29.17488962739817),
    @(11, 13.44479732892777, 8.115758399545598, 8.494474389501338, 38.0695386733229, 3.6896801525485, 28.66252743321432, 10.78838076944836, 10.45943112592059, 11.71254849829986, 15.45762155752524, 29.12148539062632),
    @(12, 13.50694270833723, 8.127207154844569, 8.499892633132307, 38.05149621953588, 3.689246760445915, 28.63983816295829, 10.78346990499494, 10.50533736828677, 11.71252712889711, 15.46991310059421, 29.10200314048614),
    @(13, 13.49357266039737, 8.124744678190988, 8.498719342609585, 38.05533471970309, 3.689339724125845, 28.64469537155176, 10.78452308762445, 10.4954606102253, 11.71252338845258, 15.46725529026391, 29.10616603527599),
    @(14, 13.4499163282644, 8.116701741104773, 8.494917274896137, 38.06803301241096, 3.689644328049718, 28.66064755569058, 10.78797473934952, 10.46321229734591, 11.7125429793024, 15.45862800147891, 29.11986772427227),
    @(15, 13.42313529103668, 8.111765839365285, 8.49260711342006, 38.07594946701459, 3.689832005665669, 28.67050461764266, 10.79010204190836, 10.44343073720601, 11.71257942559807, 15.45337470237742, 29.12835689427293),
    @(16, 13.26908813333391, 8.083338810290076, 8.479657467848087, 38.12343582628587, 3.690924383516939, 28.72830795816769, 10.80249333281975, 10.32966275217267, 11.71316374942026, 15.42375249030267, 29.17848334544665),
    @(17, 13.17411713227662, 8.065776054321876, 8.471973378253118, 38.15445840054463, 3.691609610923308, 28.76494325235914, 10.81027430826645, 10.25953845062073, 11.71385806791419, 15.40601560274506, 29.21055302017566),
    @(18, 13.11933877631088, 8.055629938905289, 8.467650158685142, 38.17299644358632, 3.692009293719437, 28.78644667870355, 10.81481575369674, 10.21909617693474, 11.71438116591752, 15.3959757672882, 29.22948306162905),
    @(19, 13.10076717329328, 8.052187091169825, 8.46620304092168, 38.17939238360364, 3.692145575291808, 28.79380153035822, 10.81636476333804, 10.20538572440538, 11.71457956076583, 15.39260444960057, 29.23597562631135),
    @(20, 13.18424327030819, 8.067650265468775, 8.472781398493352, 38.15108408979948, 3.691536092344875, 28.76099867652206, 10.80943917890532, 10.26701484434022, 11.71377135653461, 15.40788700671819, 29.20708900879205),
    @(21, 13.46274774229997, 8.119066104890349, 8.496030139298872, 38.0642743754292, 3.689554629551617, 28.65594411357159, 10.78695818332325, 10.4726904257501, 11.71253213226215, 15.4611555619252, 29.11582309654629),
    @(22, 13.64301317130349, 8.152252545379088, 8.512064558356116, 38.01373205043986, 3.688308848792273, 28.59112913607564, 10.77285072746606, 10.60586906263024, 11.71281725603844, 15.49736978823223, 29.06049363776793),
    @(23, 13.54698071824671, 8.134579450621954, 8.50343076797275, 38.04014050369486, 3.688969255047631, 28.62537038708241, 10.78032673591386, 10.53491575019192, 11.71256522933994, 15.47791553943512, 29.08962869553209),
    @(24, 13.17966579418746, 8.066803087180437, 8.472415798096252, 38.15260742674273, 3.691569312222331, 28.76278064459688, 10.80981652886038, 10.26363515883256, 11.71381017266508, 15.4070404544342, 29.20865355346937),
    @(25, 12.77783822113218, 7.99200677314699, 8.442659054267095, 38.30047425420076, 3.694586565618506, 28.92750827413558, 10.84415246683529, 9.967039589939967, 11.71988753689704, 15.33680420976069, 29.35555843935219)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
